$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.1666666666666667
$ws.Range("M2").Value = 0.6666666666666666
$ws.Range("Q2").Value = 0.3333333333333333
$ws.Range("W2").Value = 1

$ws.Range("F3").Value = 0.2
$ws.Range("M3").Value = 0.4
$ws.Range("Q3").Value = 0.4
$ws.Range("W3").Value = 0.6

$ws.Range("F4").Value = 0.1818181818181818
$ws.Range("M4").Value = 0.5
$ws.Range("Q4").Value = 0.3636363636363636
$ws.Range("W4").Value = 0.7499999999999999

$ws.Range("F5").Value = 0.1923076923076923
$ws.Range("M5").Value = 0.4347826086956522
$ws.Range("Q5").Value = 0.3846153846153846
$ws.Range("W5").Value = 0.6521739130434783

$ws.Range("F6").Value = 0.06572774036705124
$ws.Range("M6").Value = 0.7120505206430552
$ws.Range("Q6").Value = 0.142410104128611
$ws.Range("W6").Value = 0.7973386012536817
